# Auto-generated edit script applying numeric corrections to Ravana_Profits leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 5159.6
$ws.Range("I2").Value = 2400
$ws.Range("J2").Value = 6999.3335
$ws.Range("K2").Value = 2400
$ws.Range("L2").Value = 6999.3335
$ws.Range("M2").Value = -2287
$ws.Range("N2").Value = -7225.3335

# row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

# row 12
$ws.Range("H12").Value = 588.6
$ws.Range("I12").Value = 397
$ws.Range("J12").Value = 876
$ws.Range("K12").Value = 397
$ws.Range("L12").Value = 876
$ws.Range("M12").Value = -227
$ws.Range("N12").Value = -1216

# row 32
$ws.Range("H32").Value = 8855.143
$ws.Range("J32").Value = 8997.200000000001
$ws.Range("L32").Value = 8997.200000000001
$ws.Range("N32").Value = -9649.200000000001

# row 64
$ws.Range("H64").Value = 3992.5
$ws.Range("J64").Value = 3991.4285
$ws.Range("L64").Value = 3991.4285
$ws.Range("N64").Value = -4487.4285

# row 67
$ws.Range("H67").Value = 3992.5
$ws.Range("J67").Value = 3991.4285
$ws.Range("L67").Value = 3991.4285
$ws.Range("N67").Value = -5707.4285

# row 92
$ws.Range("H92").Value = 365.33334
$ws.Range("I92").Value = 356.81818
$ws.Range("J92").Value = 388.75
$ws.Range("K92").Value = 356.81818
$ws.Range("L92").Value = 388.75
$ws.Range("M92").Value = 891.18182
$ws.Range("N92").Value = -2884.75

# row 121
$ws.Range("H121").Value = 1964.3334
$ws.Range("J121").Value = 1964.3334
$ws.Range("L121").Value = 5893.0002
$ws.Range("N121").Value = -9387.0002

# row 135
$ws.Range("H135").Value = 1865.6666
$ws.Range("I135").Value = 1865.6666
$ws.Range("K135").Value = 16790.9994
$ws.Range("M135").Value = -14255.9994

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 972.9286
$ws.Range("I2").Value = 1028.8462
$ws.Range("J2").Value = 246
$ws.Range("K2").Value = 1028.8462
$ws.Range("L2").Value = 246
$ws.Range("M2").Value = -915.8462
$ws.Range("N2").Value = -472

# row 32
$ws.Range("H32").Value = 7736.8
$ws.Range("I32").Value = 7736.8
$ws.Range("K32").Value = 7736.8
$ws.Range("M32").Value = -7449.8

# row 45
$ws.Range("H45").Value = 1849.8
$ws.Range("I45").Value = 1849.8
$ws.Range("K45").Value = 1849.8
$ws.Range("M45").Value = -1472.8

# row 61
$ws.Range("H61").Value = 2050.3125
$ws.Range("I61").Value = 1431.25
$ws.Range("K61").Value = 1431.25
$ws.Range("M61").Value = -1219.25

# row 74
$ws.Range("H74").Value = 22215888
$ws.Range("I74").Value = 24992398
$ws.Range("K74").Value = 24992398
$ws.Range("M74").Value = -24991524

# row 77
$ws.Range("H77").Value = 22215888
$ws.Range("I77").Value = 24992398
$ws.Range("K77").Value = 124961990
$ws.Range("M77").Value = -124957622

# row 88
$ws.Range("H88").Value = 2000
$ws.Range("J88").Value = 2000
$ws.Range("L88").Value = 2000
$ws.Range("N88").Value = -2812

# row 91
$ws.Range("H91").Value = 2000
$ws.Range("J91").Value = 2000
$ws.Range("L91").Value = 2000
$ws.Range("N91").Value = -4808

# row 110
$ws.Range("H110").Value = 860.1429000000001
$ws.Range("I110").Value = 836.8333
$ws.Range("K110").Value = 836.8333
$ws.Range("M110").Value = 1208.1667

# row 116
$ws.Range("H116").Value = 972.9286
$ws.Range("I116").Value = 1028.8462
$ws.Range("J116").Value = 246
$ws.Range("K116").Value = 1028.8462
$ws.Range("L116").Value = 246
$ws.Range("M116").Value = 1265.1538
$ws.Range("N116").Value = -4834

# row 133
$ws.Range("H133").Value = 84749.25
$ws.Range("J133").Value = 84749.25
$ws.Range("L133").Value = 84749.25
$ws.Range("N133").Value = -89809.25

# row 136
$ws.Range("H136").Value = 2050.3125
$ws.Range("I136").Value = 1431.25
$ws.Range("K136").Value = 4293.75
$ws.Range("M136").Value = -1743.75

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 972.9286
$ws.Range("I3").Value = 1028.8462
$ws.Range("J3").Value = 246
$ws.Range("K3").Value = 1028.8462
$ws.Range("L3").Value = 246
$ws.Range("M3").Value = -914.8462
$ws.Range("N3").Value = -474

$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# row 55
$ws.Range("H55").Value = 24982
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 24982
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 24982
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -25612

# row 86
$ws.Range("H86").Value = 25444.428
$ws.Range("I86").Value = 9169.714
$ws.Range("K86").Value = 9169.714
$ws.Range("M86").Value = -8046.714

# row 89
$ws.Range("H89").Value = 25444.428
$ws.Range("I89").Value = 9169.714
$ws.Range("K89").Value = 45848.57
$ws.Range("M89").Value = -40232.57

# row 99
$ws.Range("H99").Value = 6000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 6000
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -8996

# row 126
$ws.Range("H126").Value = 6000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 18000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -22940

# row 133
$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -130060

$ws = $wb.Worksheets.Item("CUL")
# row 107
$ws.Range("H107").Value = 86.625
$ws.Range("I107").Value = 49
$ws.Range("J107").Value = 92
$ws.Range("K107").Value = 147
$ws.Range("L107").Value = 276
$ws.Range("M107").Value = 1773
$ws.Range("N107").Value = -4116

# row 131
$ws.Range("H131").Value = 1131.375
$ws.Range("J131").Value = 1832.6666
$ws.Range("L131").Value = 5497.9998
$ws.Range("N131").Value = -15577.9998

# row 132
$ws.Range("H132").Value = 5211.75
$ws.Range("J132").Value = 6249.25
$ws.Range("L132").Value = 56243.25
$ws.Range("N132").Value = -61303.25

$ws = $wb.Worksheets.Item("GSM")
# row 136
$ws.Range("H136").Value = 51424.9
$ws.Range("J136").Value = 51424.9
$ws.Range("L136").Value = 154274.7
$ws.Range("N136").Value = -159374.7

$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 2598
$ws.Range("I61").Value = 1997.75
$ws.Range("K61").Value = 1997.75
$ws.Range("M61").Value = -1795.75

# row 113
$ws.Range("H113").Value = 2598
$ws.Range("I113").Value = 1997.75
$ws.Range("K113").Value = 1997.75
$ws.Range("M113").Value = 172.25

# row 136
$ws.Range("H136").Value = 8891655
$ws.Range("I136").Value = 8891655
$ws.Range("K136").Value = 26674965
$ws.Range("M136").Value = -26672415
